# Auto-generated Excel COM-interop edit script
# Updates the cryptos list (prices / 1h volume %) per the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / percentage / label / link cells -----------------------
# These values are not valid bare numbers (e.g. contain "%", extra dots,
# letters, or surrounding spaces) so Excel keeps them as text automatically.
$plainUpdates = @{
    D2 = '69.699.43'
    E2 = '  -1.65%  '
    D3 = '3.516.67'
    E3 = '  -1.39%  '
    E4 = '  -0.08%  '
    E5 = '  +5.51%  '
    E6 = '  +1.36%  '
    E7 = '  +0.54%  '
    E8 = '  -0.09%  '
    E9 = '  -3.30%  '
    E10 = '  +0.65%  '
    E11 = '  -2.28%  '
    E12 = '  -2.87%  '
    E13 = '  +0.36%  '
    D14 = '4.079.95'
    E14 = '  -1.26%  '
    E15 = '  +8.04%  '
    D16 = '69.756.23'
    E16 = '  -1.53%  '
    E17 = '  -1.06%  '
    E18 = '  -1.28%  '
    D19 = '3.501.26'
    E19 = '  -3.30%  '
    E20 = '  -0.28%  '
    E21 = '  -1.67%  '
    B22 = 'Litecoin'
    C22 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    E22 = '  +13.98%  '
    B23 = 'InternetComputer(DFINITY)'
    C23 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    E23 = '  -2.79%  '
    E24 = '  +2.82%  '
    E25 = '  +3.20%  '
    E26 = '  +5.14%  '
    E27 = '  -2.02%  '
    E28 = '  +4.48%  '
    E29 = '  +4.34%  '
    E30 = '  -3.29%  '
    E31 = '  +1.36%  '
    E32 = '  +3.13%  '
    E33 = '  -0.11%  '
    E34 = '  -0.35%  '
    E35 = '  -5.57%  '
    D37 = '3.663.79'
    E37 = '  +0.94%  '
    B38 = 'Stacks'
    C38 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    E38 = '  +5.85%  '
    B39 = 'Bittensor'
    C39 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    E39 = '  -2.91%  '
    E40 = '  -4.16%  '
    D41 = '0.0₃0789'
    E41 = '  -1.91%  '
    E42 = '  -4.51%  '
    E43 = '  -0.41%  '
    E44 = '  -0.58%  '
    E45 = '  -0.61%  '
    E46 = '  +3.16%  '
    E47 = '  -3.46%  '
    E48 = '  -5.89%  '
    E49 = '  +0.47%  '
    E50 = '  -1.74%  '
    B51 = 'FLOKI'
    C51 = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
    E51 = '  -5.66%  '
}
foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# --- Numeric-looking price cells -----------------------------------------
# These source values (e.g. "617.27") would otherwise be auto-converted to
# a real Number by Excel on assignment. The sheet stores Price as text, so
# force a Text number format while writing, then restore the default style
# (so no permanent formatting change is left behind on the cell).
$numericTextUpdates = @{
    D5 = '617.27'
    D10 = '0.656'
    D11 = '53.44'
    D13 = '9.56'
    D15 = '619.75'
    D17 = '19.02'
    D18 = '12.60'
    D22 = '107.48'
    D23 = '17.15'
    D26 = '3.09'
    D27 = '11.00'
    D28 = '9.72'
    D29 = '34.15'
    D32 = '3.93'
    D34 = '63.45'
    D38 = '3.65'
    D39 = '515.33'
    D44 = '0.0467'
    D47 = '3.36'
    D50 = '132.09'
    D51 = '0.000241'
}
foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
    $cell.Style = "Normal"
}

"Updated " + ($plainUpdates.Count + $numericTextUpdates.Count) + " cells"
